# Applies two changes captured in the commit:
#  1. Slide 16's table (the "Total ..." summary table) is switched to a
#     different built-in PowerPoint table style.
#  2. The deck's theme colour scheme is swapped from the "Integral" theme
#     to the "Office Theme" palette (the other theme already bundled with
#     this file, previously only used by the notes master).

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 -------------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{BB1A62FA-E9AE-4151-BFBF-88CBBEC09E97}")

# --- 2. Swap the theme colour scheme ("Integral" -> "Office Theme") ------
$themeColors = $p.Slides.Item(1).ThemeColorScheme
$themeColors.Colors(1).RGB  = 0        # dk1      #000000
$themeColors.Colors(2).RGB  = 16777215 # lt1      #FFFFFF
$themeColors.Colors(3).RGB  = 6968388  # dk2      #44546A
$themeColors.Colors(4).RGB  = 15132391 # lt2      #E7E6E6
$themeColors.Colors(5).RGB  = 13998939 # accent1  #5B9BD5
$themeColors.Colors(6).RGB  = 3243501  # accent2  #ED7D31
$themeColors.Colors(7).RGB  = 10855845 # accent3  #A5A5A5
$themeColors.Colors(8).RGB  = 49407    # accent4  #FFC000
$themeColors.Colors(9).RGB  = 12874308 # accent5  #4472C4
$themeColors.Colors(10).RGB = 4697456  # accent6  #70AD47
$themeColors.Colors(11).RGB = 12673797 # hlink    #0563C1
$themeColors.Colors(12).RGB = 7491477  # folHlink #954F72
